$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new row at 70 (J-RAN-12), shifting old rows 70-142 down
#    to 71-143.
# ---------------------------------------------------------------------
$ws.Rows.Item(70).Insert()

# Renumber column A for the rows that just shifted down by one (old
# A-values are "frozen" text, they don't auto-increment on row insert).
for ($r = 71; $r -le 143; $r++) {
    $ws.Cells.Item($r, 1).Value = $r
}

# ---------------------------------------------------------------------
# 2. Populate the newly inserted row 70.
# ---------------------------------------------------------------------
$ws.Cells.Item(70, 1).Value = 70
$ws.Cells.Item(70, 2).Value = "J-RAN-12"
$ws.Cells.Item(70, 3).Value = -81.8184739
$ws.Cells.Item(70, 4).Value = 7.2680162
$ws.Cells.Item(70, 5).Value = "Jicaron"
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 8).Value = 0

# ---------------------------------------------------------------------
# 3. Append the brand-new row 144 (JIC-STREAM-CAMP-NO-T-03).
# ---------------------------------------------------------------------
$ws.Cells.Item(144, 1).Value = 144
$ws.Cells.Item(144, 2).Value = "JIC-STREAM-CAMP-NO-T-03"
$ws.Cells.Item(144, 5).Value = "Jicaron"
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 1
$ws.Cells.Item(144, 4).Value = 7.286715
$ws.Cells.Item(144, 3).Value = -81.802422

# C/D of the new rows are text-formatted columns -- match that here too.
$ws.Cells.Item(144, 3).NumberFormat = "@"
$ws.Cells.Item(144, 4).NumberFormat = "@"

# ---------------------------------------------------------------------
# 4. Misc view / print housekeeping to mirror the authored edit.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 110
$ws.Range("H144").Select()

$ws.PageSetup.Orientation = 1
